# The commit removes the stray "Test"/"Kid" smoke-test submission
# (Submission RefID 1186076, dated 26/02/2017 | 14:25) from row 2 of the
# tuckshop export, which shifts every subsequent order up by one row and
# shrinks the used range from A1:S6 to A1:S5. It also updates the saved
# window/selection state: the sheet is left zoomed to 125% with the new
# top data row (row 2) selected as a whole row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the test order in row 2; Excel shifts rows 3-6 up to 2-5 and
# prunes the now-unused shared strings automatically.
$ws.Rows(2).Delete()

# Leave the new row 2 (first real order) selected as a full row, and the
# view zoomed to 125%, matching the saved sheet state after the edit.
$ws.Rows(2).Select()
$excel.ActiveWindow.Zoom = 125
